# Applies a row-content rearrangement to the active worksheet:
#  - Rows 9 and 10 swap their full contents (columns A..AY)
#  - Rows 33, 34, 35, 36 are cyclically rotated so that each row takes on
#    the content that used to belong to the row below it (33<-34, 34<-35,
#    35<-36, 36<-33)
#
# Cell-by-cell Value2 read/write is used because bulk Range.Value array
# transfer is not reliable in this host; Value2 correctly round-trips
# both numeric and text cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1    # column A
$lastCol  = 51   # column AY

# Columns holding a plain "YYYY-MM-DD" text value (Startdatum / Slutdatum).
# These are identical across every row pair touched below, so writing to
# them is skipped to avoid Excel auto-coercing the text into a date serial
# number when it gets re-assigned through Value2.
$dateTextCols = @(25, 27)   # Y, AA

function Get-RowValues($ws, $rowIndex, $firstCol, $lastCol) {
    $values = @()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $values += $ws.Cells.Item($rowIndex, $col).Value2
    }
    return $values
}

function Set-RowValues($ws, $rowIndex, $firstCol, $lastCol, $values, $skipCols) {
    $i = 0
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        if ($skipCols -notcontains $col) {
            $ws.Cells.Item($rowIndex, $col).Value2 = $values[$i]
        }
        $i++
    }
}

# --- Swap rows 9 and 10 ---
$row9  = Get-RowValues $ws 9  $firstCol $lastCol
$row10 = Get-RowValues $ws 10 $firstCol $lastCol

Set-RowValues $ws 9  $firstCol $lastCol $row10 $dateTextCols
Set-RowValues $ws 10 $firstCol $lastCol $row9  $dateTextCols

# --- Rotate rows 33, 34, 35, 36 (33<-34, 34<-35, 35<-36, 36<-33) ---
$row33 = Get-RowValues $ws 33 $firstCol $lastCol
$row34 = Get-RowValues $ws 34 $firstCol $lastCol
$row35 = Get-RowValues $ws 35 $firstCol $lastCol
$row36 = Get-RowValues $ws 36 $firstCol $lastCol

Set-RowValues $ws 33 $firstCol $lastCol $row34 $dateTextCols
Set-RowValues $ws 34 $firstCol $lastCol $row35 $dateTextCols
Set-RowValues $ws 35 $firstCol $lastCol $row36 $dateTextCols
Set-RowValues $ws 36 $firstCol $lastCol $row33 $dateTextCols
